$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Rename the inline picture in each header/footer story.
# Footers.Item(1)/Headers.Item(1) = the "default" (primary) story,
# Footers.Item(2)/Headers.Item(2) = the "first page" story.
# The picture's Name must be selected (InlineShape.Select -> Selection.InlineShapes)
# for the rename to actually persist into the saved package.

function Rename-InlinePicture($range, [string]$newName) {
    $ishp = $range.InlineShapes.Item(1)
    $ishp.Select()
    $sel = $word.Selection
    $sel.InlineShapes.Item(1).Name = $newName
}

# Footers: PearsonLogo pictures, image1.png -> image2.png
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# Headers: BTec_Logo-Orange pictures, image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers.Item(1).Range "image1.jpg"
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"
